$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 116, shifting existing rows 116..168 down to 118..170
$ws.Rows.Item(116).Resize(2).Insert()

# Fill in the two newly inserted rows (116 and 117) with a new week of data
# Row 116: "Primera" quality
$ws.Range("A116").Value = 5
$ws.Range("B116").Value = "Macroferia Regional de Talca"
$ws.Range("C116").Value = "Maule"
$ws.Range("D116").Value = 44460
$ws.Range("D116").NumberFormat = $ws.Range("D118").NumberFormat
$ws.Range("E116").Value = 7
$ws.Range("F116").Value = 100112006
$ws.Range("G116").Value = "Repollo"
$ws.Range("H116").Value = "Crespo record"
$ws.Range("I116").Value = "Primera"
$ws.Range("J116").Value = 3000
$ws.Range("K116").Value = 500
$ws.Range("L116").Value = 500
$ws.Range("M116").Value = 500
$ws.Range("N116").Value = "$/unidad"
$ws.Range("O116").Value = "Región del Maule"
$ws.Range("P116").Value = 500
$ws.Range("Q116").Value = 1
$ws.Range("R116").Value = "Hortaliza"

# Row 117: "Segunda" quality
$ws.Range("A117").Value = 5
$ws.Range("B117").Value = "Macroferia Regional de Talca"
$ws.Range("C117").Value = "Maule"
$ws.Range("D117").Value = 44460
$ws.Range("D117").NumberFormat = $ws.Range("D119").NumberFormat
$ws.Range("E117").Value = 7
$ws.Range("F117").Value = 100112006
$ws.Range("G117").Value = "Repollo"
$ws.Range("H117").Value = "Crespo record"
$ws.Range("I117").Value = "Segunda"
$ws.Range("J117").Value = 3000
$ws.Range("K117").Value = 300
$ws.Range("L117").Value = 300
$ws.Range("M117").Value = 300
$ws.Range("N117").Value = "$/unidad"
$ws.Range("O117").Value = "Región del Maule"
$ws.Range("P117").Value = 300
$ws.Range("Q117").Value = 1
$ws.Range("R117").Value = "Hortaliza"
